# Actualización automática 2025-10-06 15:30:20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 14 ("SAL SOLUBLE") is removed; the TOTAL row shifts up from row 15 to row 14
# and its PRESUPUESTO/VENTA/POR CUMPLIR/CUMPLIMIENTO figures are recalculated.
$ws.Rows.Item(14).Delete()

# Column width tweaks (D: 14->13, E: 24->22, F: 24->25).
# ColumnWidth is expressed in characters and Excel pads it by 5/6 of a
# character for the default font before persisting to the sheet XML, so we
# back that padding out to land on the exact target widths.
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 24.166666666666668

# Refreshed VENTA (D) / POR CUMPLIR (E) / CUMPLIMIENTO (F) figures.
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 5504.61890386263
$ws.Range("F3").Value = 0

$ws.Range("D4").Value = -143.74
$ws.Range("E4").Value = 1447.7686065816
$ws.Range("F4").Value = -0.1102276432238724

$ws.Range("D5").Value = -86.41
$ws.Range("E5").Value = 236.41
$ws.Range("F5").Value = -0.5760666666666666

$ws.Range("D6").Value = 394.37
$ws.Range("E6").Value = 455.47419682004
$ws.Range("F6").Value = 0.4640497652106818

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 709.368813030059
$ws.Range("F7").Value = 0

$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 388.107983534392
$ws.Range("F10").Value = 0

$ws.Range("D11").Value = 790.78
$ws.Range("E11").Value = 2715.88949822329
$ws.Range("F11").Value = 0.2255074224704271

$ws.Range("C12").Value = 38542.25
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 38542.25
$ws.Range("F12").Value = 0

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1332.52398144409
$ws.Range("F13").Value = 0

# Row 14 is now the TOTAL row (previously row 15); refresh its totals.
$ws.Range("C14").Value = 54483.76774946896
$ws.Range("D14").Value = 955
$ws.Range("E14").Value = 53528.76774946896
$ws.Range("F14").Value = 0.01752815635642798
